$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "점령당한 학교에서 도망쳐 나온 조야는 자신의 아버지의 동료와 마주치게 되고, 좋지 않은 소식을 전해 듣게 된다.`n"
$ws.Range("D3").Value = "그건 바로, 조야의 어머니는 생사가 불투명하고, 아버지는 학생들을 구하고자 이미 페테르헤임 고등학교로 향했다는 소식이었다.`n"
